# Brought model up to version 2.1.1
# - Renamed "SYSoCCtaSC" sheet to "SYSoCCtaSC-electricity"
# - Added a new "SYSoCCtaSC-buildings" sheet (copy of the electricity sheet,
#   repurposed for distributed solar / retrofitting labor cost share)
# - Added a new title row to the "About" sheet describing the new variable

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename the existing "SYSoCCtaSC" sheet to "SYSoCCtaSC-electricity"
# ------------------------------------------------------------------
$elecSheet = $wb.Worksheets.Item("SYSoCCtaSC")
$elecSheet.Name = "SYSoCCtaSC-electricity"

# ------------------------------------------------------------------
# 2. Create the new "SYSoCCtaSC-buildings" sheet by copying the
#    electricity sheet (keeps formatting/styles) and then stripping
#    out the content that doesn't apply, replacing it with the new
#    "distributed solar" row.
# ------------------------------------------------------------------
$elecSheet.Copy([System.Reflection.Missing]::Value, $elecSheet)
$bldgSheet = $wb.Worksheets.Item($elecSheet.Index + 1)
$bldgSheet.Name = "SYSoCCtaSC-buildings"

# Clear out the old rows (values only, formatting/styles are retained)
$bldgSheet.Range("A3:A14").ClearContents()
$bldgSheet.Range("B3:B5").ClearContents()
$bldgSheet.Range("B6:B7").ClearContents()
$bldgSheet.Range("B8:B13").ClearContents()
$bldgSheet.Range("B14").ClearContents()
# The old rows 15-17 (crude oil / heavy fuel oil / municipal solid waste
# pass-through formulas) are no longer relevant on this sheet
$bldgSheet.Rows("15:17").Delete()

# New row 2: "distributed solar" with its labor cost share formula
$bldgSheet.Range("A2").Value = "distributed solar"
$bldgSheet.Range("B2").Formula = "=AVERAGE(Data!B64:B65)"

# Header row height was tightened now that the header text is shorter
$bldgSheet.Rows(1).RowHeight = 45

# ------------------------------------------------------------------
# 3. Insert a new title row on the "About" sheet describing the new
#    "distributed solar" variable, shifting all subsequent rows down.
# ------------------------------------------------------------------
$aboutSheet = $wb.Worksheets.Item("About")

# Capture existing hyperlinks so they can be re-anchored after the
# row insert (Excel's Insert() shifts cell values but not the
# hyperlink anchors automatically).
$existingLinks = @()
foreach ($h in $aboutSheet.Hyperlinks) {
    $existingLinks += ,@($h.Range.Row, $h.Range.Column, $h.Address)
}

$aboutSheet.Rows("2:2").Insert()
$aboutSheet.Range("A2").Value = "SYSoCCtaSC Share of Distributed Solar and Retrofitting Costs that is Labor"

$aboutSheet.Hyperlinks.Delete()
foreach ($link in $existingLinks) {
    $newRow = [int]$link[0] + 1
    $col = [int]$link[1]
    $target = $aboutSheet.Cells.Item($newRow, $col)
    $aboutSheet.Hyperlinks.Add($target, $link[2]) | Out-Null
    # Hyperlinks.Add() re-applies its own built-in hyperlink style; restore
    # the workbook's named "Hyperlink" cell style so formatting matches
    # the original (non-auto) style used throughout this sheet.
    $target.Style = "Hyperlink"
}
